# Insert a new column M (pushing the existing "path" column to N) and
# populate it with a relative-path header + two sample relative paths.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank column at M; existing column M ("path") and its
# data move to column N automatically, carrying their formatting along.
$ws.Columns.Item(13).Insert() | Out-Null

# Data rows (2 & 3) of the freshly inserted column should carry the
# default/normal style (no explicit style), matching the rest of the
# un-styled data cells.
$ws.Range("M2").Style = "Normal"
$ws.Range("M3").Style = "Normal"

# Populate the new column's contents.
$ws.Range("M1").Value = "./foo/bar/test"
$ws.Range("M2").Value = "relative property path 1"
$ws.Range("M3").Value = "relative property path 2"

# Restore/update the selection to reflect the newly added last column.
$ws.Range("N1:N1048576").Select() | Out-Null
